# Project Sample Project is saved. Rules sheet: row 11's Rule-name cell (B11)
# is changed from the shared string "R40" to the text "1".
# A leading apostrophe forces Excel to store the literal "1" as text
# (rather than coercing it to the number 1), which is what the target
# workbook shows: B11 keeps t="s" (shared string) pointing at a new
# <si> entry whose text is "1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
